$d = $word.ActiveDocument

function Rename-LogoShape($InlineShape, $OldName, $NewName) {
    # InlineShape.Name only ever round-trips the <wp:docPr name="..."/>
    # attribute; the nested <pic:cNvPr name="..."/> is only kept in sync
    # with the real (floating) Shape object's Name. Converting to a Shape
    # and back lets us update both copies of the name in one go, exactly
    # as Word itself does when an image is renamed.
    $floating = $InlineShape.ConvertToShape()
    if ($floating.Name -eq $OldName) {
        $floating.Name = $NewName
    }
    [void]$floating.ConvertToInlineShape()
}

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            Rename-LogoShape $shp "image1.jpg" "image2.jpg"
        }
    }
    foreach ($ftr in $sec.Footers) {
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            Rename-LogoShape $shp "image2.png" "image1.png"
        }
    }
}
